# Update stats for 2026-02 (row 27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6566
$ws.Range("D27").Value = 6131568
$ws.Range("E27").Value = 933.8361254949741
$ws.Range("F27").Value = 10.35294117647059
$ws.Range("H27").Value = 25.76267452800294
